$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the leading block of paragraphs:
#    Heading1 "Manipulacion de texto csv o txt"
#    Heading2 "Abrir archivos con csv"
#    empty paragraph
#    "nombre edad ciudad"
#    "Juan 28 Madrid"
#    "Ana 22 Barcelona"
# ---------------------------------------------------------------------------
$startP = $d.Paragraphs(1).Range.Start
$endP = $d.Paragraphs(6).Range.End
$d.Range($startP, $endP).Delete()

# ---------------------------------------------------------------------------
# 2) Split "para saltar tabulaciones agrega el parametro delimiter al
#    pd.read asi: pd.read_csv(source, delimiter='\t')" with a manual line
#    break before "pd.read_csv", adding extra leading spaces.
#    (Range.Text assignment is used instead of Find's ReplaceWith so that
#    straight quotes are not auto-corrected into curly quotes.)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("para saltar tabulaciones agrega el parametro delimiter al pd.read asi: pd.read_csv(source, delimiter='\t')") | Out-Null
$rng.Text = "para saltar tabulaciones agrega el parametro delimiter al pd.read asi:`v               pd.read_csv(source, delimiter='\t')"

# ---------------------------------------------------------------------------
# 3) Split "Para archivos con comentarios con el signo #, se debe saltar las
#    filas con el parametro skiprows = n" with a manual line break.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Para archivos con comentarios con el signo #, se debe saltar las filas con el parametro skiprows = n") | Out-Null
$rng.Text = "Para archivos con comentarios con el signo #, `v              se debe saltar las filas con el parametro skiprows = n"

# ---------------------------------------------------------------------------
# 4) Delete the "Manipular archivos excel" section entirely, from its
#    Heading2 through the final "Para guardar un dataframe..." paragraph
#    (right before "Descripcion preliminar de los datos").
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -eq "Manipular archivos excel`r") { $startPara = $i }
    if ($txt -eq "Para guardar un dataframe en un excel se convierte los datos en un dataframe y se utiliza la propiedad .to_excel, donde la propiendad index se deja en false si no se quiere agregar los numeros a las filas`r") { $endPara = $i }
}
$s = $d.Paragraphs($startPara).Range.Start
$e = $d.Paragraphs($endPara).Range.End
$d.Range($s, $e).Delete()

# ---------------------------------------------------------------------------
# 5) Split "Si quero obtener los valores unicos de una columna utilizo el
#    metodo unique, especificando la columna, ejemplo: df["nombre"].unique():"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('Si quero obtener los valores unicos de una columna utilizo el metodo unique, especificando la columna, ejemplo: df["nombre"].unique():') | Out-Null
$rng.Text = "Si quero obtener los valores unicos de una columna utilizo el metodo unique, e`v              specificando la columna, ejemplo: df[`"nombre`"].unique():"

# ---------------------------------------------------------------------------
# 6) Split "Para contar la cantidad de datos faltantes por columna utilizar
#    los metodos .isnull().sum():"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Para contar la cantidad de datos faltantes por columna utilizar los metodos .isnull().sum():") | Out-Null
$rng.Text = "Para contar la cantidad de datos faltantes por columna `v              utilizar los metodos .isnull().sum():"

# ---------------------------------------------------------------------------
# 7) Split "Si quiero eliminar las columnas con  datos faltantes se utiliza
#    el parametro axis=1 en el metodo .dropna()"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Si quiero eliminar las columnas con  datos faltantes se utiliza el parametro axis=1 en el metodo .dropna()") | Out-Null
$rng.Text = "Si quiero eliminar las columnas con  datos faltantes se utiliza el parametro axis=1`v               en el metodo .dropna()"

# ---------------------------------------------------------------------------
# 8) Split "Se puede rellenar los datos faltantes utilizando el metodo
#    fillna() y especificar que deseas rellenar"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Se puede rellenar los datos faltantes utilizando el metodo fillna() y especificar que deseas rellenar") | Out-Null
$rng.Text = "Se puede rellenar los datos faltantes utilizando el metodo fillna()`v               y especificar que deseas rellenar"

# ---------------------------------------------------------------------------
# 9) Split "Para especificar los datos en donde se deban eliminar las
#    columnas se utiliza la propiedad subset, ejemplo de subset= [Nombre]:"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Para especificar los datos en donde se deban eliminar las columnas se utiliza la propiedad subset, ejemplo de subset= [Nombre]:") | Out-Null
$rng.Text = "Para especificar los datos en donde se deban eliminar las columnas`v                            se utiliza la propiedad subset, ejemplo de subset= [Nombre]:"

Write-Output "done"
